# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Overview "Status" column text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US" (shared across sheets via the shared
#    string table).
#  - Each language sheet (zh-cn / de-de) gets its "Latest Target File" /
#    "Latest Handback File" / "Latest Handback DateTime" columns populated
#    for both data rows, with hyperlinks added on the target-file cells.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6188d9f95941f127f5d845ece6da8fbcf8b3ce64/e2e"

# --- Overview sheet: flip the handoff status text ------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).AutoFit() | Out-Null
$overview.Columns.Item(6).AutoFit() | Out-Null

# --- Language sheets: zh-cn and de-de --------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; Handback = "2016-08-30 22:40:39" },
    @{ Name = "de-de"; Handback = "2016-08-30 22:40:46" }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)
    $xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2." + $lang.Name + ".xlf"

    # Rebuild the hyperlink collection so new entries land in worksheet
    # order (row by row): A2, I2, A3, I3.
    $ws.Hyperlinks.Delete()

    $ws.Range("A2").Value = "a.md"
    $ws.Hyperlinks.Add($ws.Range("A2"), $baseUrl + "/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null

    $ws.Range("I2").Value = "a.md"
    $ws.Hyperlinks.Add($ws.Range("I2"), $baseUrl + "/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
    $ws.Range("J2").Value = $xlf
    $ws.Range("K2").Value = $lang.Handback

    $ws.Range("A3").Value = "b.md"
    $ws.Hyperlinks.Add($ws.Range("A3"), $baseUrl + "/b.md", [Type]::Missing, [Type]::Missing, "b.md") | Out-Null

    $ws.Range("I3").Value = "a.md"
    $ws.Hyperlinks.Add($ws.Range("I3"), $baseUrl + "/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
    $ws.Range("J3").Value = $xlf
    $ws.Range("K3").Value = $lang.Handback

    $ws.Columns.Item(3).AutoFit() | Out-Null
    $ws.Columns.Item(10).AutoFit() | Out-Null
}

Write-Host "Handback report generated."
